# Update ML "Label" (column C) values on Sheet2 of the BTC_vader_50 calc workbook
# after adjusting the upstream ML labels. Excel recalculates the dependent
# D/E/H/I/J formula columns automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$changes = @{
    2  = 0
    5  = 1
    8  = 1
    11 = 1
    12 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    25 = 1
    29 = 1
    30 = 1
    32 = 1
    33 = 1
    34 = 1
    37 = 1
    38 = 1
    40 = 1
    41 = 0
    45 = 1
    47 = 0
    49 = 0
    56 = 0
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 3).Value = $changes[$row]
}

$excel.CalculateFull()
